# Generate Report for Handoff
# Update the localization-status report: rows that were "low" priority and
# awaiting handoff have now been handed off, so bump their Priority to "ht"
# and refresh the "Latest Handoff Datetime" timestamp for each locale sheet.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 (Priority column E, Latest Handoff Datetime column H)
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-16 08:31:18"

# de-de sheet: rows 4-7 (Priority column E, Latest Handoff Datetime column H)
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-16 08:31:23"
